$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 629-684 ---
# The whole Primera/Segunda record block shifts down by one record: a brand-new
# record is inserted at rows 629-630, every later record keeps its row-pair slot
# position but now holds the data that used to sit one record above it, and the
# record that used to be last (rows 683-684) is appended as new rows 685-686 below.
$updates = @(
    @(629, 4, 45106),
    @(629, 10, 1800),
    @(630, 4, 45106),
    @(630, 10, 1400),
    @(631, 4, 44545),
    @(631, 10, 2500),
    @(632, 4, 44545),
    @(632, 10, 1460),
    @(633, 4, 44307),
    @(633, 10, 3400),
    @(633, 11, 550),
    @(633, 13, 575),
    @(633, 16, 288),
    @(634, 4, 44307),
    @(634, 10, 1600),
    @(634, 11, 450),
    @(634, 12, 500),
    @(634, 13, 475),
    @(634, 16, 238),
    @(635, 4, 45068),
    @(635, 10, 2000),
    @(635, 11, 500),
    @(635, 12, 600),
    @(635, 13, 550),
    @(635, 16, 275),
    @(636, 4, 45068),
    @(636, 10, 1400),
    @(636, 11, 400),
    @(636, 12, 450),
    @(636, 13, 425),
    @(636, 16, 212),
    @(637, 4, 44286),
    @(637, 10, 3400),
    @(638, 4, 44286),
    @(639, 4, 44267),
    @(639, 10, 3000),
    @(639, 11, 450),
    @(639, 12, 500),
    @(639, 13, 475),
    @(639, 16, 238),
    @(640, 4, 44267),
    @(640, 10, 1600),
    @(640, 11, 350),
    @(640, 12, 400),
    @(640, 13, 375),
    @(640, 16, 188),
    @(641, 4, 44533),
    @(641, 10, 3200),
    @(641, 11, 550),
    @(641, 12, 600),
    @(641, 13, 575),
    @(641, 16, 288),
    @(642, 4, 44533),
    @(642, 10, 1500),
    @(642, 11, 450),
    @(642, 12, 500),
    @(642, 13, 475),
    @(642, 16, 238),
    @(643, 4, 44284),
    @(643, 10, 2700),
    @(643, 11, 450),
    @(643, 12, 500),
    @(643, 13, 475),
    @(643, 16, 238),
    @(644, 4, 44284),
    @(644, 10, 1440),
    @(644, 11, 350),
    @(644, 12, 400),
    @(644, 13, 375),
    @(644, 16, 188),
    @(645, 4, 44305),
    @(645, 10, 2800),
    @(645, 11, 550),
    @(645, 12, 600),
    @(645, 13, 575),
    @(645, 16, 288),
    @(646, 4, 44305),
    @(646, 10, 1480),
    @(646, 11, 450),
    @(646, 12, 500),
    @(646, 13, 475),
    @(646, 16, 238),
    @(647, 4, 44249),
    @(647, 10, 3200),
    @(648, 4, 44249),
    @(649, 4, 44442),
    @(649, 10, 3400),
    @(650, 4, 44442),
    @(650, 10, 1600),
    @(651, 4, 44559),
    @(651, 10, 2500),
    @(652, 4, 44559),
    @(652, 10, 1500),
    @(653, 4, 44291),
    @(653, 10, 2800),
    @(653, 11, 450),
    @(653, 12, 500),
    @(653, 13, 475),
    @(653, 16, 238),
    @(654, 4, 44291),
    @(654, 10, 1460),
    @(654, 11, 350),
    @(654, 12, 400),
    @(654, 13, 375),
    @(654, 16, 188),
    @(655, 4, 44704),
    @(655, 10, 2540),
    @(655, 11, 600),
    @(655, 13, 650),
    @(655, 16, 325),
    @(656, 4, 44704),
    @(656, 10, 1400),
    @(656, 11, 500),
    @(656, 12, 550),
    @(656, 13, 525),
    @(656, 16, 262),
    @(657, 4, 44846),
    @(657, 10, 2400),
    @(657, 11, 650),
    @(657, 13, 675),
    @(657, 16, 338),
    @(658, 4, 44846),
    @(658, 10, 1360),
    @(658, 11, 550),
    @(658, 12, 600),
    @(658, 13, 575),
    @(658, 16, 288),
    @(659, 4, 44769),
    @(659, 10, 2500),
    @(659, 11, 600),
    @(659, 12, 700),
    @(659, 13, 650),
    @(659, 16, 325),
    @(660, 4, 44769),
    @(660, 10, 1400),
    @(660, 11, 500),
    @(660, 12, 550),
    @(660, 13, 525),
    @(660, 16, 262),
    @(661, 4, 45104),
    @(661, 10, 1600),
    @(662, 4, 45104),
    @(662, 10, 1200),
    @(663, 4, 45096),
    @(663, 10, 2000),
    @(663, 11, 550),
    @(663, 12, 600),
    @(663, 13, 575),
    @(663, 16, 288),
    @(664, 4, 45096),
    @(664, 10, 1400),
    @(664, 11, 450),
    @(664, 12, 500),
    @(664, 13, 475),
    @(664, 16, 238),
    @(665, 4, 44272),
    @(665, 10, 3400),
    @(665, 11, 450),
    @(665, 12, 500),
    @(665, 13, 475),
    @(665, 16, 238),
    @(666, 4, 44272),
    @(666, 10, 1600),
    @(666, 11, 350),
    @(666, 12, 400),
    @(666, 13, 375),
    @(666, 16, 188),
    @(667, 4, 44725),
    @(667, 10, 2560),
    @(667, 11, 600),
    @(667, 12, 700),
    @(667, 13, 650),
    @(667, 16, 325),
    @(668, 4, 44725),
    @(668, 10, 1400),
    @(668, 11, 500),
    @(668, 12, 550),
    @(668, 13, 525),
    @(668, 16, 262),
    @(669, 4, 44643),
    @(669, 10, 2460),
    @(669, 11, 500),
    @(669, 12, 600),
    @(669, 13, 550),
    @(669, 16, 275),
    @(670, 4, 44643),
    @(670, 10, 1260),
    @(670, 11, 400),
    @(670, 12, 450),
    @(670, 13, 425),
    @(670, 16, 212),
    @(671, 4, 44449),
    @(671, 10, 3500),
    @(671, 11, 450),
    @(671, 12, 500),
    @(671, 13, 475),
    @(671, 16, 238),
    @(672, 4, 44449),
    @(672, 10, 1600),
    @(672, 11, 350),
    @(672, 12, 400),
    @(672, 13, 375),
    @(672, 16, 188),
    @(673, 4, 44837),
    @(673, 10, 2400),
    @(673, 11, 650),
    @(673, 12, 700),
    @(673, 13, 675),
    @(673, 16, 338),
    @(674, 4, 44837),
    @(674, 10, 1400),
    @(674, 11, 550),
    @(674, 12, 600),
    @(674, 13, 575),
    @(674, 16, 288),
    @(675, 4, 44648),
    @(675, 10, 2500),
    @(676, 4, 44648),
    @(676, 10, 1340),
    @(677, 4, 45054),
    @(677, 10, 2000),
    @(678, 4, 45054),
    @(679, 4, 45049),
    @(679, 10, 2060),
    @(679, 11, 500),
    @(679, 12, 600),
    @(679, 13, 550),
    @(679, 16, 275),
    @(680, 4, 45049),
    @(680, 10, 1400),
    @(680, 11, 400),
    @(680, 12, 450),
    @(680, 13, 425),
    @(680, 16, 212),
    @(681, 4, 44881),
    @(681, 10, 2100),
    @(681, 11, 650),
    @(681, 12, 700),
    @(681, 13, 675),
    @(681, 16, 338),
    @(682, 4, 44881),
    @(682, 10, 1460),
    @(682, 11, 550),
    @(682, 12, 600),
    @(682, 13, 575),
    @(682, 16, 288),
    @(683, 4, 44991),
    @(683, 10, 2000),
    @(684, 4, 44991),
    @(684, 10, 1500),
    @(685, 1, 8),
    @(685, 2, "Terminal La Palmera de La Serena"),
    @(685, 3, "Coquimbo"),
    @(685, 4, 44641),
    @(685, 5, 4),
    @(685, 6, 100112009),
    @(685, 7, "Acelga"),
    @(685, 8, "Sin especificar"),
    @(685, 9, "Primera"),
    @(685, 10, 2500),
    @(685, 11, 500),
    @(685, 12, 600),
    @(685, 13, 550),
    @(685, 14, "`$/atado 1,5 a 2 kilos"),
    @(685, 15, "Provincia del Elquí"),
    @(685, 16, 275),
    @(685, 17, 2),
    @(685, 18, "Hortaliza"),
    @(686, 1, 8),
    @(686, 2, "Terminal La Palmera de La Serena"),
    @(686, 3, "Coquimbo"),
    @(686, 4, 44641),
    @(686, 5, 4),
    @(686, 6, 100112009),
    @(686, 7, "Acelga"),
    @(686, 8, "Sin especificar"),
    @(686, 9, "Segunda"),
    @(686, 10, 1300),
    @(686, 11, 400),
    @(686, 12, 450),
    @(686, 13, 425),
    @(686, 14, "`$/atado 1,5 a 2 kilos"),
    @(686, 15, "Provincia del Elquí"),
    @(686, 16, 212),
    @(686, 17, 2),
    @(686, 18, "Hortaliza")
)

foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}

# Column D holds dates; the two brand-new rows need the same date/time number format
# already used by every other cell in column D so they display and round-trip as dates.
$ws.Cells.Item(685, 4).NumberFormat = $ws.Cells.Item(684, 4).NumberFormat
$ws.Cells.Item(686, 4).NumberFormat = $ws.Cells.Item(684, 4).NumberFormat
